$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("VTStFES")

$ws.Range("B6").Formula = "=B2"
$ws.Range("C6").Formula = "=C2"
$ws.Range("D6").Formula = "=D2"
$ws.Range("E6").Formula = "=E2"
$ws.Range("F6").Formula = "=F2"
$ws.Range("G6").Formula = "=G2"
$ws.Range("H6").Formula = "=H2"

[void]$ws.Range("B6:H6").Select()

$about = $wb.Worksheets.Item("About")
[void]$about.Range("K6").Select()
